# Update average_county_temperature (K) and derived worst/best WHP COP
# (R, S) values for rows 6-8 using refreshed NOAA temperature data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K6").Value = 19.36574074074073
$ws.Range("R6").Value = 2.499849249524808
$ws.Range("S6").Value = 2.870901733221348

$ws.Range("K7").Value = 1.925925925925943
$ws.Range("R7").Value = 2.004846509671994
$ws.Range("S7").Value = 2.229613377609108

$ws.Range("K8").Value = 1.925925925925943
